$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

Set-TextValue $ws.Range('D2') '29.580.69'
Set-TextValue $ws.Range('E2') '  -1.12%  '
Set-TextValue $ws.Range('D3') '1.858.04'
Set-TextValue $ws.Range('E3') '  -0.28%  '
Set-TextValue $ws.Range('D4') '0.9986'
Set-TextValue $ws.Range('E4') '  -0.19%  '
Set-TextValue $ws.Range('D5') '242.43'
Set-TextValue $ws.Range('E5') '  -1.03%  '
Set-TextValue $ws.Range('E6') '  -3.76%  '
Set-TextValue $ws.Range('D7') '0.9995'
Set-TextValue $ws.Range('E7') '  -0.10%  '
Set-TextValue $ws.Range('D8') '0.07603'
Set-TextValue $ws.Range('E8') '  +0.07%  '
Set-TextValue $ws.Range('D9') '0.2995'
Set-TextValue $ws.Range('E9') '  -0.34%  '
Set-TextValue $ws.Range('D10') '24.67'
Set-TextValue $ws.Range('E10') '  -0.57%  '
Set-TextValue $ws.Range('D11') '0.07739'
Set-TextValue $ws.Range('E11') '  +0.95%  '
Set-TextValue $ws.Range('D12') '1.860.12'
Set-TextValue $ws.Range('E12') '  -0.31%  '
Set-TextValue $ws.Range('D13') '0.6948'
Set-TextValue $ws.Range('E13') '  -0.08%  '
Set-TextValue $ws.Range('D14') '5.035'
Set-TextValue $ws.Range('E14') '  -0.98%  '
Set-TextValue $ws.Range('D15') '83.73'
Set-TextValue $ws.Range('E15') '  -0.43%  '
Set-TextValue $ws.Range('D16') '0.00001002'
Set-TextValue $ws.Range('E16') '  +2.46%  '
Set-TextValue $ws.Range('D17') '2.113.17'
Set-TextValue $ws.Range('E17') '  -0.52%  '
Set-TextValue $ws.Range('D18') '6.277'
Set-TextValue $ws.Range('E18') '  +1.93%  '
Set-TextValue $ws.Range('D19') '29.606.77'
Set-TextValue $ws.Range('E19') '  -1.12%  '
Set-TextValue $ws.Range('D20') '234.85'
Set-TextValue $ws.Range('E20') '  -0.88%  '
Set-TextValue $ws.Range('D21') '12.59'
Set-TextValue $ws.Range('E21') '  -1.07%  '
Set-TextValue $ws.Range('D22') '0.9997'
Set-TextValue $ws.Range('E22') '  -0.10%  '
Set-TextValue $ws.Range('D23') '7.677'
Set-TextValue $ws.Range('E23') '  -1.05%  '
Set-TextValue $ws.Range('D24') '0.9994'
Set-TextValue $ws.Range('D25') '155.95'
Set-TextValue $ws.Range('E25') '  -1.95%  '
Set-TextValue $ws.Range('E26') '  -3.44%  '
Set-TextValue $ws.Range('D27') '8.494'
Set-TextValue $ws.Range('E27') '  -1.49%  '
Set-TextValue $ws.Range('D28') '17.78'
Set-TextValue $ws.Range('E28') '  -1.15%  '
Set-TextValue $ws.Range('D29') '1.479'
Set-TextValue $ws.Range('E29') '  -1.22%  '
Set-TextValue $ws.Range('D30') '0.05829'
Set-TextValue $ws.Range('E30') '  -4.00%  '
Set-TextValue $ws.Range('D31') '1.263'
Set-TextValue $ws.Range('E31') '  -1.87%  '
Set-TextValue $ws.Range('D32') '4.140'
Set-TextValue $ws.Range('D33') '4.036'
Set-TextValue $ws.Range('E33') '  -2.01%  '
Set-TextValue $ws.Range('D34') '1.908'
Set-TextValue $ws.Range('E34') '  +1.13%  '
Set-TextValue $ws.Range('D35') '1.172'
Set-TextValue $ws.Range('E35') '  -1.10%  '
Set-TextValue $ws.Range('D36') '0.7224'
Set-TextValue $ws.Range('E36') '  -2.19%  '
Set-TextValue $ws.Range('D37') '2.587'
Set-TextValue $ws.Range('E37') '  -1.05%  '
Set-TextValue $ws.Range('D38') '1.250.15'
Set-TextValue $ws.Range('E38') '  +2.72%  '
Set-TextValue $ws.Range('D39') '2.808'
Set-TextValue $ws.Range('E39') '  -0.49%  '
Set-TextValue $ws.Range('E40') '  +0.45%  '
Set-TextValue $ws.Range('D41') '0.9099'
Set-TextValue $ws.Range('E41') '  -0.85%  '
Set-TextValue $ws.Range('D42') '6.138'
Set-TextValue $ws.Range('E42') '  -4.03%  '
Set-TextValue $ws.Range('D43') '0.9990'
Set-TextValue $ws.Range('E43') '  -0.21%  '
Set-TextValue $ws.Range('D44') '2.022.16'
Set-TextValue $ws.Range('E44') '  -0.57%  '
Set-TextValue $ws.Range('D45') '68.26'
Set-TextValue $ws.Range('E45') '  +0.75%  '
Set-TextValue $ws.Range('D46') '101.52'
Set-TextValue $ws.Range('E46') '  -0.25%  '
Set-TextValue $ws.Range('D47') '7.375'
Set-TextValue $ws.Range('E47') '  -5.28%  '
Set-TextValue $ws.Range('D48') '0.4063'
Set-TextValue $ws.Range('E48') '  -0.73%  '
Set-TextValue $ws.Range('D49') '9.194'
Set-TextValue $ws.Range('E49') '  -0.62%  '

# Row 50/51 swap: RenderToken <-> BabyDogeCoin
Set-TextValue $ws.Range("B50") 'BabyDogeCoin'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D50") '0.00000000117'
Set-TextValue $ws.Range("E50") '  -2.97%  '

Set-TextValue $ws.Range("B51") 'RenderToken'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D51") '1.715'
Set-TextValue $ws.Range("E51") '  +1.59%  '
